$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 560776.4399999999
$ws.Range("J2").Value = 2115.2856
$ws.Range("L2").Value = 2115.2856
$ws.Range("N2").Value = -2341.2856

# Row 33
$ws.Range("H33").Value = 3473.5938
$ws.Range("I33").Value = 3900.2593
$ws.Range("J33").Value = 1169.6
$ws.Range("K33").Value = 3900.2593
$ws.Range("L33").Value = 1169.6
$ws.Range("M33").Value = -3671.2593
$ws.Range("N33").Value = -1627.6

# Row 42
$ws.Range("H42").Value = 1348
$ws.Range("I42").Value = 865
$ws.Range("K42").Value = 2595
$ws.Range("M42").Value = -2365

# Row 81
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()

# Row 86
$ws.Range("H86").Value = 126811.875
$ws.Range("I86").Value = 168349.33
$ws.Range("K86").Value = 168349.33
$ws.Range("M86").Value = -167226.33

# Row 89
$ws.Range("H89").Value = 126811.875
$ws.Range("I89").Value = 168349.33
$ws.Range("K89").Value = 841746.6499999999
$ws.Range("M89").Value = -836130.6499999999

# Row 92
$ws.Range("H92").Value = 250824.75
$ws.Range("I92").Value = 500399.5
$ws.Range("K92").Value = 500399.5
$ws.Range("M92").Value = -499151.5

# Row 96
$ws.Range("H96").Value = 1939
$ws.Range("J96").Value = 2325
$ws.Range("L96").Value = 6975
$ws.Range("N96").Value = -9721

# Row 107
$ws.Range("H107").Value = 3209
$ws.Range("J107").Value = 3444
$ws.Range("L107").Value = 3444
$ws.Range("N107").Value = -7284

# Row 113
$ws.Range("H113").Value = 4248.273
$ws.Range("I113").Value = 3890.3845
$ws.Range("J113").Value = 4765.222
$ws.Range("K113").Value = 3890.3845
$ws.Range("L113").Value = 4765.222
$ws.Range("M113").Value = -636.3845000000001
$ws.Range("N113").Value = -11273.222

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 19959.15
$ws.Range("I32").Value = 19938.357
$ws.Range("K32").Value = 19938.357
$ws.Range("M32").Value = -19651.357

# Row 61
$ws.Range("H61").Value = 14079.895
$ws.Range("I61").Value = 15819.9375
$ws.Range("K61").Value = 15819.9375
$ws.Range("M61").Value = -15607.9375

# Row 74
$ws.Range("H74").Value = 27976.46
$ws.Range("I74").Value = 27976.46
$ws.Range("K74").Value = 27976.46
$ws.Range("M74").Value = -27102.46

# Row 77
$ws.Range("H77").Value = 27976.46
$ws.Range("I77").Value = 27976.46
$ws.Range("K77").Value = 139882.3
$ws.Range("M77").Value = -135514.3

# Row 107
$ws.Range("H107").Value = 200000
$ws.Range("J107").Value = 200000
$ws.Range("L107").Value = 200000
$ws.Range("N107").Value = -207680

# Row 110
$ws.Range("H110").Value = 2000.32
$ws.Range("I110").Value = 1422.2106
$ws.Range("K110").Value = 1422.2106
$ws.Range("M110").Value = 622.7893999999999

# Row 112
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("N112").Value = 0
$ws.Range("L112").ClearContents()

# Row 132
$ws.Range("H132").Value = 34567
$ws.Range("I132").Value = 39103.57
$ws.Range("K132").Value = 117310.71
$ws.Range("M132").Value = -114780.71

# Row 136
$ws.Range("H136").Value = 14079.895
$ws.Range("I136").Value = 15819.9375
$ws.Range("K136").Value = 47459.8125
$ws.Range("M136").Value = -44909.8125

$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 2809.8667
$ws.Range("I64").Value = 2010.1428
$ws.Range("J64").Value = 3509.625
$ws.Range("K64").Value = 2010.1428
$ws.Range("L64").Value = 3509.625
$ws.Range("M64").Value = -1785.1428
$ws.Range("N64").Value = -3959.625

# Row 67
$ws.Range("H67").Value = 2809.8667
$ws.Range("I67").Value = 2010.1428
$ws.Range("J67").Value = 3509.625
$ws.Range("K67").Value = 2010.1428
$ws.Range("L67").Value = 3509.625
$ws.Range("M67").Value = -1230.1428
$ws.Range("N67").Value = -5069.625

# Row 99
$ws.Range("H99").Value = 2495.4119
$ws.Range("I99").Value = 2740.111
$ws.Range("J99").Value = 2220.125
$ws.Range("K99").Value = 2740.111
$ws.Range("L99").Value = 2220.125
$ws.Range("M99").Value = -1242.111
$ws.Range("N99").Value = -5216.125

# Row 112
$ws.Range("H112").Value = 135954.17
$ws.Range("J112").Value = 135954.17
$ws.Range("L112").Value = 135954.17
$ws.Range("N112").Value = -138908.17

# Row 134
$ws.Range("H134").Value = 2744.1072
$ws.Range("I134").Value = 2493.44
$ws.Range("J134").Value = 4833
$ws.Range("K134").Value = 7480.32
$ws.Range("L134").Value = 14499
$ws.Range("M134").Value = -4945.32
$ws.Range("N134").Value = -19569

$ws = $wb.Worksheets.Item("CRP")
# Row 11
$ws.Range("H11").Value = 50
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

# Row 31
$ws.Range("H31").Value = 2290.4138
$ws.Range("I31").Value = 1834.1482
$ws.Range("K31").Value = 1834.1482
$ws.Range("M31").Value = -1539.1482

# Row 34
$ws.Range("H34").Value = 2290.4138
$ws.Range("I34").Value = 1834.1482
$ws.Range("K34").Value = 1834.1482
$ws.Range("M34").Value = -1632.1482

# Row 132
$ws.Range("H132").Value = 2800.1072
$ws.Range("J132").Value = 3712.8572
$ws.Range("L132").Value = 11138.5716
$ws.Range("N132").Value = -16198.5716

$ws = $wb.Worksheets.Item("CUL")
# Row 124
$ws.Range("H124").Value = 8125
$ws.Range("J124").Value = 10000
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3842.8572
$ws.Range("I80").Value = 1750
$ws.Range("K80").Value = 1750
$ws.Range("M80").Value = -752

# Row 83
$ws.Range("H83").Value = 3842.8572
$ws.Range("I83").Value = 1750
$ws.Range("K83").Value = 8750
$ws.Range("M83").Value = -3758

# Row 97
$ws.Range("H97").Value = 924.6667
$ws.Range("I97").Value = 779.3889
$ws.Range("K97").Value = 779.3889
$ws.Range("M97").Value = -283.3889

# Row 107
$ws.Range("H107").Value = 28788.195
$ws.Range("I107").Value = 46200.773
$ws.Range("J107").Value = 1425.5714
$ws.Range("K107").Value = 46200.773
$ws.Range("L107").Value = 1425.5714
$ws.Range("M107").Value = -44280.773
$ws.Range("N107").Value = -5265.5714

# Row 113
$ws.Range("H113").Value = 156665.31
$ws.Range("I113").Value = 93786.27
$ws.Range("K113").Value = 93786.27
$ws.Range("M113").Value = -91616.27

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 2556.2727
$ws.Range("I61").Value = 1623.4642
$ws.Range("K61").Value = 1623.4642
$ws.Range("M61").Value = -1421.4642

# Row 113
$ws.Range("H113").Value = 2556.2727
$ws.Range("I113").Value = 1623.4642
$ws.Range("K113").Value = 1623.4642
$ws.Range("M113").Value = 546.5358000000001

# Row 132
$ws.Range("H132").Value = 18114.977
$ws.Range("I132").Value = 21140.797
$ws.Range("K132").Value = 63422.391
$ws.Range("M132").Value = -60892.391

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 3749.8438
$ws.Range("I136").Value = 3320.6365
$ws.Range("K136").Value = 9961.9095
$ws.Range("M136").Value = -7411.9095
